$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly cryptos price/volume refresh.
# Column D ("Price") values are free-form text (dotted-thousands, variable
# decimal places) -- force text format so Excel does not reinterpret them
# as numbers (which would also strip meaningful trailing zeros).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.674.01"
$ws.Range("E2").Value = "  -2.28%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.804.59"
$ws.Range("E3").Value = "  -2.55%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.61"
$ws.Range("E5").Value = "  +0.40%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.601"
$ws.Range("E6").Value = "  -1.19%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "38.67"
$ws.Range("E8").Value = "  -7.27%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.316"
$ws.Range("E9").Value = "  +3.04%  "

$ws.Range("E10").Value = "  -1.97%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0990"
$ws.Range("E11").Value = "  -2.15%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.065.66"
$ws.Range("E12").Value = "  -2.69%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.784.21"
$ws.Range("E13").Value = "  -3.69%  "

$ws.Range("E14").Value = "  -1.85%  "

$ws.Range("E15").Value = "  -5.30%  "

$ws.Range("E16").Value = "  -3.09%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "34.658.50"
$ws.Range("E17").Value = "  -2.42%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.25"
$ws.Range("E18").Value = "  -0.81%  "

$ws.Range("E19").Value = "  -2.50%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "239.12"
$ws.Range("E20").Value = "  -3.19%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.68"
$ws.Range("E21").Value = "  -3.22%  "

$ws.Range("E22").Value = "  +0.02%  "

$ws.Range("E23").Value = "  +0.12%  "

$ws.Range("E24").Value = "  +2.09%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "172.52"
$ws.Range("E25").Value = "  +2.10%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.68"
$ws.Range("E26").Value = "  -2.78%  "

$ws.Range("E27").Value = "  -3.85%  "

$ws.Range("E28").Value = "  -2.76%  "

$ws.Range("E29").Value = "  +8.23%  "

$ws.Range("E30").Value = "  -0.08%  "

$ws.Range("E31").Value = "  +1.63%  "

$ws.Range("E32").Value = "  -0.38%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.92"
$ws.Range("E33").Value = "  -3.20%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.25"
$ws.Range("E34").Value = "  +15.09%  "

$ws.Range("E35").Value = "  -5.89%  "

$ws.Range("E36").Value = "  -1.34%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "90.49"
$ws.Range("E37").Value = "  -8.27%  "

$ws.Range("E38").Value = "  +4.31%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.304.94"
$ws.Range("E39").Value = "  -4.07%  "

$ws.Range("E40").Value = "  -2.16%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.46"
$ws.Range("E41").Value = "  -1.00%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.952"
$ws.Range("E42").Value = "  -5.04%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "14.12"
$ws.Range("E43").Value = "  -3.52%  "

$ws.Range("E44").Value = "  -9.95%  "

$ws.Range("E45").Value = "  -5.35%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.13"
$ws.Range("E46").Value = "  -1.24%  "

$ws.Range("E47").Value = "  -1.82%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.993.46"
$ws.Range("E48").Value = "  -1.40%  "

$ws.Range("E49").Value = "  +0.03%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0670"
$ws.Range("E50").Value = "  +7.88%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "98.41"
$ws.Range("E51").Value = "  -5.00%  "
